# Auto-generated edits applying scheduled market-data refresh to Famfrit_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 222.45454
$ws.Cells.Item(6, 9).Value = 244.2
$ws.Cells.Item(6, 10).Value = 5
$ws.Cells.Item(6, 11).Value = 732.5999999999999
$ws.Cells.Item(6, 12).Value = 15
$ws.Cells.Item(6, 13).Value = -620.5999999999999
$ws.Cells.Item(6, 14).Value = -239

$ws.Cells.Item(17, 8).Value = 795621.4399999999
$ws.Cells.Item(17, 10).Value = 795621.4399999999
$ws.Cells.Item(17, 12).Value = 2386864.32
$ws.Cells.Item(17, 14).Value = -2387200.32

$ws.Cells.Item(33, 8).Value = 249.3125
$ws.Cells.Item(33, 9).Value = 150.28572
$ws.Cells.Item(33, 11).Value = 150.28572
$ws.Cells.Item(33, 13).Value = 78.71428

$ws.Cells.Item(38, 8).Value = 5109
$ws.Cells.Item(38, 10).Value = 6666.6665
$ws.Cells.Item(38, 12).Value = 19999.9995
$ws.Cells.Item(38, 14).Value = -20743.9995

$ws.Cells.Item(80, 8).Value = 4666
$ws.Cells.Item(80, 9).Value = 4199.8
$ws.Cells.Item(80, 10).Value = 4899.1
$ws.Cells.Item(80, 11).Value = 12599.4
$ws.Cells.Item(80, 12).Value = 14697.3
$ws.Cells.Item(80, 13).Value = -11601.4
$ws.Cells.Item(80, 14).Value = -16693.3

$ws.Cells.Item(81, 8).Value = 90000
$ws.Cells.Item(81, 9).Value = 90000
$ws.Cells.Item(81, 11).Value = 90000
$ws.Cells.Item(81, 13).Value = -89002

$ws.Cells.Item(83, 8).Value = 4666
$ws.Cells.Item(83, 9).Value = 4199.8
$ws.Cells.Item(83, 10).Value = 4899.1
$ws.Cells.Item(83, 11).Value = 37798.2
$ws.Cells.Item(83, 12).Value = 44091.9
$ws.Cells.Item(83, 13).Value = -32806.2
$ws.Cells.Item(83, 14).Value = -54075.9

$ws.Cells.Item(84, 8).Value = 90000
$ws.Cells.Item(84, 9).Value = 90000
$ws.Cells.Item(84, 11).Value = 270000
$ws.Cells.Item(84, 13).Value = -265008

$ws.Cells.Item(127, 8).Value = 1800.8889
$ws.Cells.Item(127, 10).Value = 4108
$ws.Cells.Item(127, 12).Value = 12324
$ws.Cells.Item(127, 14).Value = -22244

$ws.Cells.Item(138, 8).Value = 7251152
$ws.Cells.Item(138, 10).Value = 10107156
$ws.Cells.Item(138, 12).Value = 30321468
$ws.Cells.Item(138, 14).Value = -30331748

$ws.Cells.Item(141, 8).Value = 2287.7
$ws.Cells.Item(141, 9).Value = 2287.7
$ws.Cells.Item(141, 11).Value = 6863.099999999999
$ws.Cells.Item(141, 13).Value = -1683.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1093.75
$ws.Cells.Item(97, 10).Value = 1024.2
$ws.Cells.Item(97, 12).Value = 1024.2
$ws.Cells.Item(97, 14).Value = -2016.2

$ws.Cells.Item(102, 8).Value = 203816.5
$ws.Cells.Item(102, 9).Value = 336163.16
$ws.Cells.Item(102, 10).Value = 5296.5
$ws.Cells.Item(102, 11).Value = 336163.16
$ws.Cells.Item(102, 12).Value = 5296.5
$ws.Cells.Item(102, 13).Value = -334541.16
$ws.Cells.Item(102, 14).Value = -8540.5

$ws.Cells.Item(122, 8).Value = 3426.4814
$ws.Cells.Item(122, 10).Value = 4022.8667
$ws.Cells.Item(122, 12).Value = 12068.6001
$ws.Cells.Item(122, 14).Value = -16968.6001

$ws.Cells.Item(123, 8).Value = 59000
$ws.Cells.Item(123, 10).Value = 59000
$ws.Cells.Item(123, 12).Value = 59000
$ws.Cells.Item(123, 14).Value = -68800

$ws.Cells.Item(132, 8).Value = 34540944
$ws.Cells.Item(132, 9).Value = 11756.695
$ws.Cells.Item(132, 11).Value = 35270.085
$ws.Cells.Item(132, 13).Value = -32740.085

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 6710.8887
$ws.Cells.Item(99, 10).Value = 6802
$ws.Cells.Item(99, 12).Value = 6802
$ws.Cells.Item(99, 14).Value = -9798

$ws.Cells.Item(105, 8).Value = 9077.385
$ws.Cells.Item(105, 9).Value = 12311.777
$ws.Cells.Item(105, 10).Value = 1800
$ws.Cells.Item(105, 11).Value = 12311.777
$ws.Cells.Item(105, 12).Value = 1800
$ws.Cells.Item(105, 13).Value = -10564.777
$ws.Cells.Item(105, 14).Value = -5294

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 12136.25
$ws.Cells.Item(99, 9).Value = 12136.25
$ws.Cells.Item(99, 11).Value = 12136.25
$ws.Cells.Item(99, 13).Value = -10638.25

$ws.Cells.Item(126, 8).Value = 12136.25
$ws.Cells.Item(126, 9).Value = 12136.25
$ws.Cells.Item(126, 11).Value = 36408.75
$ws.Cells.Item(126, 13).Value = -33938.75

$ws.Cells.Item(132, 8).Value = 57470.73
$ws.Cells.Item(132, 10).Value = 5766.1665
$ws.Cells.Item(132, 12).Value = 17298.4995
$ws.Cells.Item(132, 14).Value = -22358.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 2500775
$ws.Cells.Item(81, 9).Value = 3334100
$ws.Cells.Item(81, 11).Value = 10002300
$ws.Cells.Item(81, 13).Value = -10001177

$ws.Cells.Item(84, 8).Value = 2500775
$ws.Cells.Item(84, 9).Value = 3334100
$ws.Cells.Item(84, 11).Value = 30006900
$ws.Cells.Item(84, 13).Value = -30001284

$ws.Cells.Item(124, 8).Value = 1748.75
$ws.Cells.Item(124, 9).Value = 998.5
$ws.Cells.Item(124, 10).Value = 2499
$ws.Cells.Item(124, 11).Value = 2995.5
$ws.Cells.Item(124, 12).Value = 7497
$ws.Cells.Item(124, 13).Value = 1914.5
$ws.Cells.Item(124, 14).Value = -17317

$ws.Cells.Item(125, 8).Value = 10333.333
$ws.Cells.Item(125, 10).Value = 10333.333
$ws.Cells.Item(125, 12).Value = 30999.999
$ws.Cells.Item(125, 14).Value = -40839.999

$ws.Cells.Item(126, 8).Value = 7000
$ws.Cells.Item(126, 9).Value = 5000
$ws.Cells.Item(126, 10).Value = 9000
$ws.Cells.Item(126, 11).Value = 15000
$ws.Cells.Item(126, 12).Value = 27000
$ws.Cells.Item(126, 13).Value = -10060
$ws.Cells.Item(126, 14).Value = -36880

$ws.Cells.Item(128, 8).Value = 115923.75
$ws.Cells.Item(128, 9).Value = 115923.75
$ws.Cells.Item(128, 11).Value = 347771.25
$ws.Cells.Item(128, 13).Value = -342791.25

$ws.Cells.Item(131, 8).Value = 1686.8
$ws.Cells.Item(131, 10).Value = 1717.2778
$ws.Cells.Item(131, 12).Value = 5151.8334
$ws.Cells.Item(131, 14).Value = -15231.8334

$ws.Cells.Item(133, 8).Value = 3682.1428
$ws.Cells.Item(133, 9).Value = 3682.1428
$ws.Cells.Item(133, 11).Value = 11046.4284
$ws.Cells.Item(133, 13).Value = -5986.428400000001

$ws.Cells.Item(134, 8).Value = 7622.263
$ws.Cells.Item(134, 10).Value = 14750.375
$ws.Cells.Item(134, 12).Value = 44251.125
$ws.Cells.Item(134, 14).Value = -54391.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2901.2693
$ws.Cells.Item(102, 9).Value = 1691.4375
$ws.Cells.Item(102, 10).Value = 4837
$ws.Cells.Item(102, 11).Value = 1691.4375
$ws.Cells.Item(102, 12).Value = 4837
$ws.Cells.Item(102, 13).Value = -69.4375
$ws.Cells.Item(102, 14).Value = -8081

$ws.Cells.Item(113, 8).Value = 3591.25
$ws.Cells.Item(113, 9).Value = 797.5
$ws.Cells.Item(113, 10).Value = 4150
$ws.Cells.Item(113, 11).Value = 797.5
$ws.Cells.Item(113, 12).Value = 4150
$ws.Cells.Item(113, 13).Value = 1372.5
$ws.Cells.Item(113, 14).Value = -8490

$ws.Cells.Item(122, 8).Value = 450
$ws.Cells.Item(122, 9).Value = 450
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 1350
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 1100
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 1826.8334
$ws.Cells.Item(132, 9).Value = 1621.6875
$ws.Cells.Item(132, 10).Value = 2237.125
$ws.Cells.Item(132, 11).Value = 4865.0625
$ws.Cells.Item(132, 12).Value = 6711.375
$ws.Cells.Item(132, 13).Value = -2335.0625
$ws.Cells.Item(132, 14).Value = -11771.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3347
$ws.Cells.Item(7, 9).Value = 2550
$ws.Cells.Item(7, 11).Value = 2550
$ws.Cells.Item(7, 13).Value = -2438

$ws.Cells.Item(40, 8).Value = 3857.8667
$ws.Cells.Item(40, 9).Value = 3682.2307
$ws.Cells.Item(40, 10).Value = 4999.5
$ws.Cells.Item(40, 11).Value = 3682.2307
$ws.Cells.Item(40, 12).Value = 4999.5
$ws.Cells.Item(40, 13).Value = -3546.2307
$ws.Cells.Item(40, 14).Value = -5271.5

$ws.Cells.Item(112, 8).Value = 59947.5
$ws.Cells.Item(112, 10).Value = 59947.5
$ws.Cells.Item(112, 12).Value = 59947.5
$ws.Cells.Item(112, 14).Value = -62901.5

$ws.Cells.Item(126, 8).Value = 3347
$ws.Cells.Item(126, 9).Value = 2550
$ws.Cells.Item(126, 11).Value = 7650
$ws.Cells.Item(126, 13).Value = -5180

$ws.Cells.Item(132, 8).Value = 5148.625
$ws.Cells.Item(132, 9).Value = 4331.3335
$ws.Cells.Item(132, 10).Value = 7600.5
$ws.Cells.Item(132, 11).Value = 12994.0005
$ws.Cells.Item(132, 12).Value = 22801.5
$ws.Cells.Item(132, 13).Value = -10464.0005
$ws.Cells.Item(132, 14).Value = -27861.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 66607
$ws.Cells.Item(46, 10).Value = 66607
$ws.Cells.Item(46, 12).Value = 66607
$ws.Cells.Item(46, 14).Value = -67069

$ws.Cells.Item(81, 8).Value = 1708.7826
$ws.Cells.Item(81, 9).Value = 1678
$ws.Cells.Item(81, 10).Value = 1742.3636
$ws.Cells.Item(81, 11).Value = 3356
$ws.Cells.Item(81, 12).Value = 3484.7272
$ws.Cells.Item(81, 13).Value = -2295
$ws.Cells.Item(81, 14).Value = -5606.727199999999

$ws.Cells.Item(84, 8).Value = 1708.7826
$ws.Cells.Item(84, 9).Value = 1678
$ws.Cells.Item(84, 10).Value = 1742.3636
$ws.Cells.Item(84, 11).Value = 16780
$ws.Cells.Item(84, 12).Value = 17423.636
$ws.Cells.Item(84, 13).Value = -11476
$ws.Cells.Item(84, 14).Value = -28031.636

$ws.Cells.Item(86, 8).Value = 11137477
$ws.Cells.Item(86, 10).Value = 18999.6
$ws.Cells.Item(86, 12).Value = 18999.6
$ws.Cells.Item(86, 14).Value = -21245.6

$ws.Cells.Item(89, 8).Value = 11137477
$ws.Cells.Item(89, 10).Value = 18999.6
$ws.Cells.Item(89, 12).Value = 94998
$ws.Cells.Item(89, 14).Value = -106230

$ws.Cells.Item(110, 8).Value = 79990
$ws.Cells.Item(110, 10).Value = 79990
$ws.Cells.Item(110, 12).Value = 79990
$ws.Cells.Item(110, 14).Value = -88170

$ws.Cells.Item(134, 8).Value = 66607
$ws.Cells.Item(134, 10).Value = 66607
$ws.Cells.Item(134, 12).Value = 199821
$ws.Cells.Item(134, 14).Value = -204891
